$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression coefficient table values (kept as text so values
# such as "44.29***" retain their significance-star suffixes). A leading
# apostrophe forces Excel to store purely-numeric-looking strings (e.g.
# "0.17") as text instead of converting them to numbers; the style is then
# reset to Normal so no stray number formatting/quote-prefix style sticks
# to the cell.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("B2") "0.17"
Set-TextValue $ws.Range("B3") "-0.01"
Set-TextValue $ws.Range("B4") "-0.09"

Set-TextValue $ws.Range("C2") "44.29***"
Set-TextValue $ws.Range("C3") "2.21***"
Set-TextValue $ws.Range("C4") "0.98"

Set-TextValue $ws.Range("D2") "-0.89"
Set-TextValue $ws.Range("D3") "0.46***"
Set-TextValue $ws.Range("D4") "0.82*"
